$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Label" header in H1, matching the style of the other headers
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the Label column: 0 for Control rows, 1 for MDD rows
$labels = @{
  2 = 0; 3 = 0; 4 = 0; 5 = 0; 6 = 0; 7 = 1; 8 = 1; 9 = 1; 10 = 1; 11 = 1;
  12 = 0; 13 = 0; 14 = 0; 15 = 0; 16 = 0; 17 = 1; 18 = 1; 19 = 1; 20 = 1; 21 = 1;
}
foreach ($r in $labels.Keys) {
  $ws.Cells.Item($r, 8).Value = $labels[$r]
}

# Update refitted D/E values with new precision
$ws.Range("D4").Value = 0.4885369609516565
$ws.Range("E4").Value = 0.4885369609516565

$ws.Range("D6").Value = 0.5265685213742143
$ws.Range("E6").Value = 0.5265685213742143

$ws.Range("D7").Value = 0.5191090209208207
$ws.Range("E7").Value = 0.4808909790791793

$ws.Range("D10").Value = 0.5008863370982628
$ws.Range("E10").Value = 0.4991136629017372

$ws.Range("D11").Value = 0.5648205115450385
$ws.Range("E11").Value = 0.4351794884549615
